$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3
$ws.Range("B3").Value = 5437306.463627676

# Clear D3 (cell becomes empty)
$ws.Range("D3").ClearContents()

# Update C4
$ws.Range("C4").Value = 8208.664146488076

# Update C5
$ws.Range("C5").Value = 12316.58249218266

# Row 7: rename "Other" -> "Biogas", update D7
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 2504.480677045021

# New row 8: "Other" with D8 value, copy style from A7 to A8
$ws.Range("A8").Value = "Other"
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("D8").Value = 1742.488088247493
